$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "62.801.72"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "  -0.55%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.444.61"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "  -1.02%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.999"
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "  -0.06%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "576.06"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "  +0.38%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "144.19"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "  -1.99%  "
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "  +0.03%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.530"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "  -2.01%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "2.440.84"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "  -1.05%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.108"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "  -3.54%  "
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "  +0.46%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "5.20"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "  -0.79%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.350"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "  -1.74%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "26.55"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "  -1.82%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.0000175"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "  -2.17%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "2.871.57"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "  -1.49%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "62.647.74"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "  -0.54%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "2.437.41"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "  -1.58%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "11.13"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "  -2.92%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "7.15"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "  -2.37%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "330.28"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "  +1.00%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "4.14"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "  -0.77%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "2.04"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "  +4.28%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "1.01"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "  +0.65%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "65.72"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "  -0.31%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "634.19"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "  +1.00%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "9.14"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  +7.55%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.0₃0968"
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "  -5.86%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.998"
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "  -0.06%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.45"
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "  -3.13%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "8.09"
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  -2.08%  "
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "  -1.45%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.138"
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "  -3.56%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "5.00"
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "  -2.49%  "
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "  +0.41%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.47"
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "  -1.62%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.376"
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "  -1.81%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "18.52"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "  -1.35%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "5.28"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "  -2.71%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "146.60"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "  -0.43%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "1.75"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "  -2.12%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "42.45"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  +1.54%  "
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "  -0.01%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "2.51"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "  -3.86%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "145.10"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "  -1.93%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "3.70"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "  -0.52%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.0527"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "  -2.50%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.600"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "  -0.75%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "19.76"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "  -4.93%  "
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "  -2.31%  "
